$d = $word.ActiveDocument

# 1) Normalize the trailing double-space to a single space, and append a
#    one-character sentinel so the new run split below does not land
#    exactly at the paragraph end (which would pick up a stray w:rsidRPr).
$replaced = $d.Content.Find.Execute(
    "again.  ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "again. #", 2
)

# 2) Locate the boundary right before the sentinel and split the run there
#    (toggling Bold on/off around InsertAfter forces a clean run split).
$r = $d.Content
$r.Find.Execute("again. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.Bold = 1
$r.InsertAfter("Make sure your 7 segment display is set to “cathode”. The default is anode which is incorrect. ")

# 3) Clear the Bold toggle across the new text AND the trailing sentinel so
#    they share identical rPr again (merging them back into one run), then
#    strip the now-redundant sentinel character with a plain text replace.
$full = $d.Range($r.Start, $r.End + 1)
$full.Bold = 0

$d.Content.Find.Execute("incorrect. #", $true, $false, $false, $false, $false, $true, 1, $false, "incorrect. ", 2)
